# Add support for zipfiles.
# Restructure the AIDA Pathology Anonymization Sheet:
#  - shift the instructions/settings block from column B to column D
#  - reorder / reword the Howto steps
#  - rework the header row (row 14): add Status, rename Image file -> OrigFile,
#    drop AnonFile and the trailing "..." column, widen the bordered header band
#  - update column widths and the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Start from a clean slate for the whole used range so that no stray cells
# from the old layout survive.
# ---------------------------------------------------------------------------
$ws.Range("A1:R14").Clear()

# ---------------------------------------------------------------------------
# Title / Howto block (was column B, now column D)
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "AIDA Pathology Anonymization Sheet"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Size = 14

$ws.Range("D2").Value = "Howto:"
$ws.Range("D2").Font.Bold = $true

$ws.Range("D3").Value = "1. Put a copy of this file in a folder. Put exported .zip files or folders named after Persons in same folder. Should contain BLOCK_STAIN (eg ""A_HE"") subfolders with one .svs file each in them."

$ws.Range("D4").Value = "2. Set Prefix for your study. AnonIDs are generated sequentially from the last AnonID given based on this prefix, or from 1 if none given (eg MYPROJ-001)."

$ws.Range("D5").Value = "3. Run aida-pat-anonexcel.py <path to this file> to fill out missing AnonID and AnonFile. Anonymized image files appear in subdirectory anon."

# Row 6: "4. " (plain) + "Your data is now Pseudonymous" (bold) + rest (plain)
$d6 = "4. Your data is now Pseudonymous because keys still exist that connect AnonIDs to Persons. Please verify that everything went as expected."
$ws.Range("D6").Value = $d6
$r = $ws.Range("D6").Characters(4, 29)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D6").Characters(33, 106)
$r.Font.Size = 11
$r.Font.Name = "Calibri"

# Row 7: "5. " (plain) + "To make your data Anonymous" (bold) + ": Delete all keys..." (plain)
#        + "Person" (bold) + " and " (plain) + "OrigFile" (bold) + " cells below..." (plain)
#        + "Study parameters" (bold) + " may not contain identifiers." (plain)
$d7 = "5. To make your data Anonymous: Delete all keys associating AnonIDs to persons, including the Person and OrigFile cells below and any other identifiers that may exist. Obviously, Study parameters may not contain identifiers."
$ws.Range("D7").Value = $d7
$r = $ws.Range("D7").Characters(4, 27)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D7").Characters(31, 64)
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D7").Characters(95, 6)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D7").Characters(101, 5)
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D7").Characters(106, 8)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D7").Characters(114, 66)
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D7").Characters(180, 16)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Calibri"
$r = $ws.Range("D7").Characters(196, 29)
$r.Font.Size = 11
$r.Font.Name = "Calibri"

$ws.Range("D8").Font.Bold = $true

# ---------------------------------------------------------------------------
# Prefix / Digits / Example block (was column B/C, now column D/E)
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "Prefix:"
$ws.Range("D9").Font.Bold = $true
$ws.Range("E9").Value = "MYPROJ-"
$ws.Range("G9").Font.Bold = $true
$ws.Range("H9").Font.Bold = $true

$ws.Range("D10").Value = "Digits:"
$ws.Range("D10").Font.Bold = $true
$ws.Range("E10").Value = 3

$ws.Range("D11").Font.Bold = $true

$ws.Range("D12").Value = "Example:"
$ws.Range("D12").Font.Bold = $true
$ws.Range("E12").Formula = "=E9&TEXT(1, REPT(""0"", E10))"

# ---------------------------------------------------------------------------
# Header row 14: Status, Person, OrigFile, AnonID, Block, Stain,
#                Study parameter 1-3, then blank bordered cells out to R
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Status"
$ws.Range("B14").Value = "Person"
$ws.Range("C14").Value = "OrigFile"
$ws.Range("D14").Value = "AnonID"
$ws.Range("E14").Value = "Block"
$ws.Range("F14").Value = "Stain"
$ws.Range("G14").Value = "Study parameter 1"
$ws.Range("H14").Value = "Study parameter 2"
$ws.Range("I14").Value = "Study parameter 3…"

$ws.Range("A14:N14").Font.Bold = $true
$ws.Range("A14:R14").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 3.5
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 9.5
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 12

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("F15").Select()
